# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets.
# Both sheets contain the same rows/data and need the same updates.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    4  = 732
    9  = 439
    10 = 40
    12 = 536
    13 = 21
    14 = 277
    16 = 334
    18 = 83
    19 = 45
    20 = 40
    22 = 79
    23 = 824
    24 = 1357
    25 = 282
    26 = 294
    28 = 59
    29 = 150
    32 = 79
    33 = 194
    34 = 231
    35 = 245
    36 = 1581
    39 = 150
    40 = 557
    42 = 3289
    43 = 391
    44 = 169
    45 = 853
    46 = 52
    47 = 38
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
